$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "MSPE" column (D). This shifts the
# old D (MSPE) -> E, old E (Nugget) -> F, old F (Non-expl var of model) -> G,
# matching how the source notebook now reports an extra VAR(OK) column ahead
# of the recomputed MSPE.
$ws.Columns.Item(4).Insert()

# New header row layout: search_radius, #PCs, VAR(OK), MSPE, S_nugget, VAR(TOTAL), VAR(DATA)
$ws.Range("D1").Value = "VAR(OK)"
$ws.Range("E1").Value = "MSPE"
$ws.Range("F1").Value = "S_nugget"

# Two brand-new trailing columns need the same header formatting (border +
# bold + centered) as the rest of row 1 - copy it over from an existing
# header cell before writing the new labels.
$ws.Range("C1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "VAR(TOTAL)"
$ws.Range("H1").Value = "VAR(DATA)"

# New VAR(OK) values (column D)
$ws.Range("D2").Value = 0.6257652903833159
$ws.Range("D3").Value = 0.9124815930662299
$ws.Range("D4").Value = 1.157086074073914
$ws.Range("D5").Value = 1.292427869341214
$ws.Range("D6").Value = 1.393382669966238
$ws.Range("D7").Value = 1.47797386397048
$ws.Range("D8").Value = 1.514747761993534
$ws.Range("D9").Value = 1.539375878697728
$ws.Range("D10").Value = 1.561999672684794

# Recomputed MSPE values (column E) - replace the shifted-over originals
$ws.Range("E2").Value = 5.999797496840213
$ws.Range("E3").Value = 5.867739111075558
$ws.Range("E4").Value = 5.762763442279576
$ws.Range("E5").Value = 5.652759981456985
$ws.Range("E6").Value = 5.576654963647206
$ws.Range("E7").Value = 5.545205571062069
$ws.Range("E8").Value = 5.528228541709547
$ws.Range("E9").Value = 5.52356707496152
$ws.Range("E10").Value = 5.51917028891843

# Column F (S_nugget) keeps the old "Nugget" values untouched - nothing to do,
# Columns.Insert already shifted them from E to F.

# New VAR(TOTAL) column (G) has no data yet for these rows
$ws.Range("G2:G10").Value = ""

# New VAR(DATA) column (H) - constant across all rows
$ws.Range("H2").Value = 6.599951515112549
$ws.Range("H3").Value = 6.599951515112549
$ws.Range("H4").Value = 6.599951515112549
$ws.Range("H5").Value = 6.599951515112549
$ws.Range("H6").Value = 6.599951515112549
$ws.Range("H7").Value = 6.599951515112549
$ws.Range("H8").Value = 6.599951515112549
$ws.Range("H9").Value = 6.599951515112549
$ws.Range("H10").Value = 6.599951515112549
